$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grab a style donor cell (style index 2: shaded fill used for B/C data cells) before we touch it ---
$ws.Range("C2").Copy() | Out-Null

# Apply that same "data cell" formatting to every B/C cell we are about to populate or
# leave blank-but-styled, so the saved workbook reuses cellXf s="2" instead of minting a
# new style.
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null

# --- Existing rows: clear stale "Number Of Test Cases" counts ---
$ws.Range("B2").Value = $null
$ws.Range("B3").Value = $null

# Web_SEARCH (row 4) now has a (blank) styled count cell
# (already formatted above via PasteSpecial)

# Web_MESSAGING / Web_QUICK_MESSAGES pick up new counts
$ws.Range("B5").Value = 63
$ws.Range("B6").Value = 38

# Web_QUICK_MESSAGES Runmode flips from Y to N
$ws.Range("C6").Value = "N"

# --- New rows describing new messaging / upload / conversation suites ---
$ws.Range("A7").Value = "Web_FILE_UPLOAD_NEW_USER"
$ws.Range("B7").Value = 84
$ws.Range("C7").Value = "N"

$ws.Range("A8").Value = "Web_FILE_UPLOAD_EXISTING_USER"
$ws.Range("C8").Value = "Y"

$ws.Range("A9").Value = "Web_URGENT_MESSAGES"
$ws.Range("B9").Value = 11
$ws.Range("C9").Value = "N"

$ws.Range("A10").Value = "Web_CONVERSATIONS"
$ws.Range("C10").Value = "N"

# --- Column widths (TSID / Number Of Test Cases columns got re-sized) ---
$ws.Columns.Item(1).ColumnWidth = 27.43
$ws.Columns.Item(2).ColumnWidth = 16.75

# --- Final selection lands on the newly added last row ---
$ws.Range("A10").Select() | Out-Null
